$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.468.15'
$ws.Range("E2").Value = '  -1.98%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.177.31'
$ws.Range("E3").Value = '  -3.63%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.07'
$ws.Range("E5").Value = '  -3.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.36'
$ws.Range("E6").Value = '  -4.86%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.174.86'
$ws.Range("E8").Value = '  -3.68%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.508'
$ws.Range("E9").Value = '  -2.23%  '
$ws.Range("E10").Value = '  -5.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.27'
$ws.Range("E11").Value = '  -3.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.455'
$ws.Range("E12").Value = '  -3.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000235'
$ws.Range("E13").Value = '  -4.89%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.31'
$ws.Range("E14").Value = '  -3.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.700.61'
$ws.Range("E15").Value = '  -3.67%  '
$ws.Range("E16").Value = '  -1.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.176.79'
$ws.Range("E17").Value = '  -3.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.484.56'
$ws.Range("E18").Value = '  -2.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.52'
$ws.Range("E19").Value = '  -4.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '455.07'
$ws.Range("E20").Value = '  -5.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.94'
$ws.Range("E21").Value = '  -1.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.702'
$ws.Range("E22").Value = '  -3.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.61'
$ws.Range("E23").Value = '  -5.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.73'
$ws.Range("E24").Value = '  -0.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.23'
$ws.Range("E25").Value = '  -1.94%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.69'
$ws.Range("E27").Value = '  -3.21%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.83'
$ws.Range("E29").Value = '  -6.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.75'
$ws.Range("E30").Value = '  -4.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.02'
$ws.Range("E31").Value = '  -7.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.26'
$ws.Range("E32").Value = '  -6.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.103'
$ws.Range("E33").Value = '  -1.92%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.38'
$ws.Range("E34").Value = '  -6.12%  '
$ws.Range("E35").Value = '  -6.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.91'
$ws.Range("E36").Value = '  -1.10%  '
$ws.Range("E37").Value = '  -3.73%  '
$ws.Range("E38").Value = '  -6.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0385'
$ws.Range("E39").Value = '  -3.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.73'
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.113'
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.99'
$ws.Range("E42").Value = '  -4.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '392.41'
$ws.Range("E43").Value = '  -7.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.799.43'
$ws.Range("E44").Value = '  -8.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '36.19'
$ws.Range("E45").Value = '  +3.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.249'
$ws.Range("E46").Value = '  -5.90%  '
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.13'
$ws.Range("E47").Value = '  -2.73%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.999'
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.98'
$ws.Range("E49").Value = '  +0.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.31'
$ws.Range("E50").Value = '  -3.62%  '
$ws.Range("E51").Value = '  -3.81%  '
